$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALERTS")
$ws.Range("A6").Value = "'2026-02-06"
$ws.Range("B6").Value = "09:51:06"
$ws.Range("C6").Value = "09:00"
$ws.Range("D6").Value = "Bathroom"
$ws.Range("E6").Value = "MINIMAL"
$ws.Range("F6").Value = "MINIMAL ALERT: Bathroom occupied, no motion > 20s."

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A125").Value = "'2026-02-06"
$ws.Range("B125").Value = "09:50:27"
$ws.Range("C125").Value = "09:00"
$ws.Range("D125").Value = "Bathroom"
$ws.Range("E125").Value = "Motion Detected"
$ws.Range("F125").Value = "Active"
$ws.Range("A126").Value = "'2026-02-06"
$ws.Range("B126").Value = "09:50:28"
$ws.Range("C126").Value = "09:00"
$ws.Range("D126").Value = "Bathroom"
$ws.Range("E126").Value = "No Motion"
$ws.Range("F126").Value = "Inactive"
$ws.Range("A127").Value = "'2026-02-06"
$ws.Range("B127").Value = "09:50:29"
$ws.Range("C127").Value = "09:00"
$ws.Range("D127").Value = "Bathroom"
$ws.Range("E127").Value = "Motion Detected"
$ws.Range("F127").Value = "Active"
$ws.Range("A128").Value = "'2026-02-06"
$ws.Range("B128").Value = "09:50:37"
$ws.Range("C128").Value = "09:00"
$ws.Range("D128").Value = "Bathroom"
$ws.Range("E128").Value = "No Motion"
$ws.Range("F128").Value = "Inactive"
$ws.Range("A129").Value = "'2026-02-06"
$ws.Range("B129").Value = "09:50:42"
$ws.Range("C129").Value = "09:00"
$ws.Range("D129").Value = "Bathroom"
$ws.Range("E129").Value = "No Motion"
$ws.Range("F129").Value = "Inactive"
$ws.Range("A130").Value = "'2026-02-06"
$ws.Range("B130").Value = "09:50:44"
$ws.Range("C130").Value = "09:00"
$ws.Range("D130").Value = "Bathroom"
$ws.Range("E130").Value = "Motion Detected"
$ws.Range("F130").Value = "Active"
$ws.Range("A131").Value = "'2026-02-06"
$ws.Range("B131").Value = "09:50:51"
$ws.Range("C131").Value = "09:00"
$ws.Range("D131").Value = "Bathroom"
$ws.Range("E131").Value = "No Motion"
$ws.Range("F131").Value = "Inactive"
$ws.Range("A132").Value = "'2026-02-06"
$ws.Range("B132").Value = "09:50:56"
$ws.Range("C132").Value = "09:00"
$ws.Range("D132").Value = "Bathroom"
$ws.Range("E132").Value = "No Motion"
$ws.Range("F132").Value = "Inactive"
$ws.Range("A133").Value = "'2026-02-06"
$ws.Range("B133").Value = "09:51:01"
$ws.Range("C133").Value = "09:00"
$ws.Range("D133").Value = "Bathroom"
$ws.Range("E133").Value = "No Motion"
$ws.Range("F133").Value = "Inactive"
$ws.Range("A134").Value = "'2026-02-06"
$ws.Range("B134").Value = "09:51:07"
$ws.Range("C134").Value = "09:00"
$ws.Range("D134").Value = "Bathroom"
$ws.Range("E134").Value = "No Motion"
$ws.Range("F134").Value = "Inactive"
$ws.Range("A135").Value = "'2026-02-06"
$ws.Range("B135").Value = "09:51:11"
$ws.Range("C135").Value = "09:00"
$ws.Range("D135").Value = "Bathroom"
$ws.Range("E135").Value = "Motion Detected"
$ws.Range("F135").Value = "Active"
$ws.Range("A136").Value = "'2026-02-06"
$ws.Range("B136").Value = "09:51:18"
$ws.Range("C136").Value = "09:00"
$ws.Range("D136").Value = "Bathroom"
$ws.Range("E136").Value = "No Motion"
$ws.Range("F136").Value = "Inactive"
$ws.Range("A137").Value = "'2026-02-06"
$ws.Range("B137").Value = "09:51:23"
$ws.Range("C137").Value = "09:00"
$ws.Range("D137").Value = "Bathroom"
$ws.Range("E137").Value = "No Motion"
$ws.Range("F137").Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A57").Value = "'2026-02-06"
$ws.Range("B57").Value = "09:50:28"
$ws.Range("C57").Value = "09:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "'71.6%"
$ws.Range("F57").Value = "Active"
$ws.Range("A58").Value = "'2026-02-06"
$ws.Range("B58").Value = "09:50:30"
$ws.Range("C58").Value = "09:00"
$ws.Range("D58").Value = "Bathroom"
$ws.Range("E58").Value = "'71.6%"
$ws.Range("F58").Value = "Active"
$ws.Range("A59").Value = "'2026-02-06"
$ws.Range("B59").Value = "09:50:41"
$ws.Range("C59").Value = "09:00"
$ws.Range("D59").Value = "Bathroom"
$ws.Range("E59").Value = "'71.5%"
$ws.Range("F59").Value = "Active"
$ws.Range("A60").Value = "'2026-02-06"
$ws.Range("B60").Value = "09:50:46"
$ws.Range("C60").Value = "09:00"
$ws.Range("D60").Value = "Bathroom"
$ws.Range("E60").Value = "'88.6%"
$ws.Range("F60").Value = "Active"
$ws.Range("A61").Value = "'2026-02-06"
$ws.Range("B61").Value = "09:50:51"
$ws.Range("C61").Value = "09:00"
$ws.Range("D61").Value = "Bathroom"
$ws.Range("E61").Value = "'80.4%"
$ws.Range("F61").Value = "Active"
$ws.Range("A62").Value = "'2026-02-06"
$ws.Range("B62").Value = "09:50:56"
$ws.Range("C62").Value = "09:00"
$ws.Range("D62").Value = "Bathroom"
$ws.Range("E62").Value = "'90.3%"
$ws.Range("F62").Value = "Active"
$ws.Range("A63").Value = "'2026-02-06"
$ws.Range("B63").Value = "09:51:01"
$ws.Range("C63").Value = "09:00"
$ws.Range("D63").Value = "Bathroom"
$ws.Range("E63").Value = "'81.4%"
$ws.Range("F63").Value = "Active"
$ws.Range("A64").Value = "'2026-02-06"
$ws.Range("B64").Value = "09:51:06"
$ws.Range("C64").Value = "09:00"
$ws.Range("D64").Value = "Bathroom"
$ws.Range("E64").Value = "'88.9%"
$ws.Range("F64").Value = "Active"
$ws.Range("A65").Value = "'2026-02-06"
$ws.Range("B65").Value = "09:51:11"
$ws.Range("C65").Value = "09:00"
$ws.Range("D65").Value = "Bathroom"
$ws.Range("E65").Value = "'85.3%"
$ws.Range("F65").Value = "Active"
$ws.Range("A66").Value = "'2026-02-06"
$ws.Range("B66").Value = "09:51:16"
$ws.Range("C66").Value = "09:00"
$ws.Range("D66").Value = "Bathroom"
$ws.Range("E66").Value = "'86.7%"
$ws.Range("F66").Value = "Active"
$ws.Range("A67").Value = "'2026-02-06"
$ws.Range("B67").Value = "09:51:21"
$ws.Range("C67").Value = "09:00"
$ws.Range("D67").Value = "Bathroom"
$ws.Range("E67").Value = "'80.4%"
$ws.Range("F67").Value = "Active"
$ws.Range("A68").Value = "'2026-02-06"
$ws.Range("B68").Value = "09:51:26"
$ws.Range("C68").Value = "09:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "'77.1%"
$ws.Range("F68").Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A57").Value = "'2026-02-06"
$ws.Range("B57").Value = "09:50:28"
$ws.Range("C57").Value = "09:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "27.5C"
$ws.Range("F57").Value = "Active"
$ws.Range("A58").Value = "'2026-02-06"
$ws.Range("B58").Value = "09:50:31"
$ws.Range("C58").Value = "09:00"
$ws.Range("D58").Value = "Bathroom"
$ws.Range("E58").Value = "27.5C"
$ws.Range("F58").Value = "Active"
$ws.Range("A59").Value = "'2026-02-06"
$ws.Range("B59").Value = "09:50:41"
$ws.Range("C59").Value = "09:00"
$ws.Range("D59").Value = "Bathroom"
$ws.Range("E59").Value = "27.5C"
$ws.Range("F59").Value = "Active"
$ws.Range("A60").Value = "'2026-02-06"
$ws.Range("B60").Value = "09:50:46"
$ws.Range("C60").Value = "09:00"
$ws.Range("D60").Value = "Bathroom"
$ws.Range("E60").Value = "27.6C"
$ws.Range("F60").Value = "Active"
$ws.Range("A61").Value = "'2026-02-06"
$ws.Range("B61").Value = "09:50:51"
$ws.Range("C61").Value = "09:00"
$ws.Range("D61").Value = "Bathroom"
$ws.Range("E61").Value = "27.5C"
$ws.Range("F61").Value = "Active"
$ws.Range("A62").Value = "'2026-02-06"
$ws.Range("B62").Value = "09:50:56"
$ws.Range("C62").Value = "09:00"
$ws.Range("D62").Value = "Bathroom"
$ws.Range("E62").Value = "27.6C"
$ws.Range("F62").Value = "Active"
$ws.Range("A63").Value = "'2026-02-06"
$ws.Range("B63").Value = "09:51:01"
$ws.Range("C63").Value = "09:00"
$ws.Range("D63").Value = "Bathroom"
$ws.Range("E63").Value = "27.6C"
$ws.Range("F63").Value = "Active"
$ws.Range("A64").Value = "'2026-02-06"
$ws.Range("B64").Value = "09:51:06"
$ws.Range("C64").Value = "09:00"
$ws.Range("D64").Value = "Bathroom"
$ws.Range("E64").Value = "27.6C"
$ws.Range("F64").Value = "Active"
$ws.Range("A65").Value = "'2026-02-06"
$ws.Range("B65").Value = "09:51:11"
$ws.Range("C65").Value = "09:00"
$ws.Range("D65").Value = "Bathroom"
$ws.Range("E65").Value = "27.7C"
$ws.Range("F65").Value = "Active"
$ws.Range("A66").Value = "'2026-02-06"
$ws.Range("B66").Value = "09:51:16"
$ws.Range("C66").Value = "09:00"
$ws.Range("D66").Value = "Bathroom"
$ws.Range("E66").Value = "27.7C"
$ws.Range("F66").Value = "Active"
$ws.Range("A67").Value = "'2026-02-06"
$ws.Range("B67").Value = "09:51:21"
$ws.Range("C67").Value = "09:00"
$ws.Range("D67").Value = "Bathroom"
$ws.Range("E67").Value = "27.7C"
$ws.Range("F67").Value = "Active"
$ws.Range("A68").Value = "'2026-02-06"
$ws.Range("B68").Value = "09:51:26"
$ws.Range("C68").Value = "09:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "27.7C"
$ws.Range("F68").Value = "Active"

$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A6").Value = "'2026-02-06"
$ws.Range("B6").Value = "09:50:27"
$ws.Range("C6").Value = "09:00"
$ws.Range("D6").Value = "Bathroom Door"
$ws.Range("E6").Value = "EXIT"
$ws.Range("F6").Value = "User EXITED Bathroom"
$ws.Range("A7").Value = "'2026-02-06"
$ws.Range("B7").Value = "09:50:31"
$ws.Range("C7").Value = "09:00"
$ws.Range("D7").Value = "Bathroom Door"
$ws.Range("E7").Value = "ENTER"
$ws.Range("F7").Value = "User ENTERED Bathroom"
